$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '29.118.65'
$ws.Range("E2").Value = '  +0.24%  '
$ws.Range("D3").Value = '1.839.31'
$ws.Range("E3").Value = '  +0.23%  '
$ws.Range("D4").Value = "'0.9978"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.14%  '
$ws.Range("D5").Value = "'243.37"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  -0.24%  '
$ws.Range("D6").Value = "'0.6257"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -1.12%  '
$ws.Range("D7").Value = "'0.9991"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.14%  '
$ws.Range("E8").Value = '  -0.95%  '
$ws.Range("D9").Value = "'0.2947"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  -0.01%  '
$ws.Range("D10").Value = "'23.37"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +2.45%  '
$ws.Range("D11").Value = "'0.07702"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  -0.46%  '
$ws.Range("D12").Value = '1.836.13'
$ws.Range("E12").Value = '  -0.24%  '
$ws.Range("D13").Value = "'5.026"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +0.59%  '
$ws.Range("D14").Value = "'0.6769"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +0.89%  '
$ws.Range("D15").Value = "'83.09"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = '  -0.19%  '
$ws.Range("D16").Value = "'0.000009374"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  -4.84%  '
$ws.Range("D17").Value = "'5.991"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  -2.09%  '
$ws.Range("D18").Value = '29.109.04'
$ws.Range("E18").Value = '  +0.08%  '
$ws.Range("D19").Value = '2.080.99'
$ws.Range("E19").Value = '  -0.31%  '
$ws.Range("D20").Value = "'12.71"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.16%  '
$ws.Range("D21").Value = "'227.42"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +0.26%  '
$ws.Range("D22").Value = "'0.9999"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +0.02%  '
$ws.Range("D23").Value = "'7.167"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  -1.01%  '
$ws.Range("D24").Value = "'0.9994"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  -0.07%  '
$ws.Range("D25").Value = "'160.24"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  -0.18%  '
$ws.Range("D26").Value = "'0.1403"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  -0.17%  '
$ws.Range("D27").Value = "'8.549"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.02%  '
$ws.Range("D28").Value = "'17.96"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  -0.01%  '
$ws.Range("E29").Value = '  -0.67%  '
$ws.Range("D30").Value = "'4.188"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +1.66%  '
$ws.Range("D31").Value = "'4.150"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +2.39%  '
$ws.Range("D32").Value = "'0.05570"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +3.70%  '
$ws.Range("D33").Value = "'1.207"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +0.45%  '
$ws.Range("D34").Value = "'0.7492"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  -0.07%  '
$ws.Range("D35").Value = "'1.850"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  -0.73%  '
$ws.Range("D36").Value = "'1.148"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.49%  '
$ws.Range("D37").Value = "'2.662"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -0.37%  '
$ws.Range("D38").Value = '1.237.34'
$ws.Range("E38").Value = '  -0.95%  '
$ws.Range("D39").Value = "'2.772"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +0.44%  '
$ws.Range("D40").Value = "'0.01786"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  -0.63%  '
$ws.Range("D41").Value = "'6.588"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  -0.28%  '
$ws.Range("D42").Value = "'0.8972"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  -0.87%  '
$ws.Range("D43").Value = "'0.9991"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  -0.15%  '
$ws.Range("D44").Value = "'102.39"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.24%  '
$ws.Range("D45").Value = '1.983.16'
$ws.Range("E45").Value = '  -0.22%  '
$ws.Range("D46").Value = "'66.58"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +2.48%  '
$ws.Range("D47").Value = "'0.00000000124"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  -0.83%  '
$ws.Range("D48").Value = "'0.5083"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  -0.56%  '
$ws.Range("D49").Value = "'0.4079"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  -0.45%  '
$ws.Range("D50").Value = "'9.088"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +0.07%  '
$ws.Range("D51").Value = "'0.05841"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +0.77%  '
